$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N4").Copy($ws.Range("O4"))
$ws.Range("O4").Value = 2022
$ws.Range("N5").Copy($ws.Range("O5"))
$ws.Range("O5").Value = 27683.5
$ws.Range("N6").Copy($ws.Range("O6"))
$ws.Range("O6").Value = 24954.3
$ws.Range("N7").Copy($ws.Range("O7"))
$ws.Range("O7").Value = 31353.200000000001
$ws.Range("N8").Copy($ws.Range("O8"))
$ws.Range("O8").Value = 20225.8
$ws.Range("N9").Copy($ws.Range("O9"))
$ws.Range("O9").Value = 21514
$ws.Range("N10").Copy($ws.Range("O10"))
$ws.Range("O10").Value = 19900.5
$ws.Range("N11").Copy($ws.Range("O11"))
$ws.Range("O11").Value = 43159.199999999997
$ws.Range("N12").Copy($ws.Range("O12"))
$ws.Range("O12").Value = 39796.5
$ws.Range("N13").Copy($ws.Range("O13"))
$ws.Range("O13").Value = 43595.8
$ws.Range("N14").Copy($ws.Range("O14"))
$ws.Range("O14").Value = 27079.3
$ws.Range("N15").Copy($ws.Range("O15"))
$ws.Range("O15").Value = 21245.1
$ws.Range("N16").Copy($ws.Range("O16"))
$ws.Range("O16").Value = 30178.6
$ws.Range("N17").Copy($ws.Range("O17"))
$ws.Range("O17").Value = 34707.199999999997
$ws.Range("N18").Copy($ws.Range("O18"))
$ws.Range("O18").Value = 31299.8
$ws.Range("N19").Copy($ws.Range("O19"))
$ws.Range("O19").Value = 35538.1
$ws.Range("N20").Copy($ws.Range("O20"))
$ws.Range("O20").Value = 21527.1
$ws.Range("N21").Copy($ws.Range("O21"))
$ws.Range("O21").Value = 20882.3
$ws.Range("N22").Copy($ws.Range("O22"))
$ws.Range("O22").Value = 21786.6
$ws.Range("N23").Copy($ws.Range("O23"))
$ws.Range("O23").Value = 26375.200000000001
$ws.Range("N24").Copy($ws.Range("O24"))
$ws.Range("O24").Value = 21717.1
$ws.Range("N25").Copy($ws.Range("O25"))
$ws.Range("O25").Value = 27335.1
$ws.Range("N26").Copy($ws.Range("O26"))
$ws.Range("O26").Value = 26244.9
$ws.Range("N27").Copy($ws.Range("O27"))
$ws.Range("O27").Value = 21435.8
$ws.Range("N28").Copy($ws.Range("O28"))
$ws.Range("O28").Value = 31143.1
$ws.Range("N29").Copy($ws.Range("O29"))
$ws.Range("O29").Value = 30942.400000000001
$ws.Range("N30").Copy($ws.Range("O30"))
$ws.Range("O30").Value = 20516.3
$ws.Range("N31").Copy($ws.Range("O31"))
$ws.Range("O31").Value = 36018
$ws.Range("N32").Copy($ws.Range("O32"))
$ws.Range("O32").Value = 21990.9
$ws.Range("N33").Copy($ws.Range("O33"))
$ws.Range("O33").Value = 19605.8
$ws.Range("N34").Copy($ws.Range("O34"))
$ws.Range("O34").Value = 24386.7
$ws.Range("N35").Copy($ws.Range("O35"))
$ws.Range("O35").Value = 51522.6
$ws.Range("N36").Copy($ws.Range("O36"))
$ws.Range("O36").Value = 42311.5
$ws.Range("N37").Copy($ws.Range("O37"))
$ws.Range("O37").Value = 56950.400000000001
$ws.Range("N38").Copy($ws.Range("O38"))
$ws.Range("O38").Value = 46649.3
$ws.Range("N39").Copy($ws.Range("O39"))
$ws.Range("O39").Value = 40968.300000000003
$ws.Range("N40").Copy($ws.Range("O40"))
$ws.Range("O40").Value = 54331.9
$ws.Range("N41").Copy($ws.Range("O41"))
$ws.Range("O41").Value = 20994.9
$ws.Range("N42").Copy($ws.Range("O42"))
$ws.Range("O42").Value = 19543.3
$ws.Range("N43").Copy($ws.Range("O43"))
$ws.Range("O43").Value = 21964.1
$ws.Range("N44").Copy($ws.Range("O44"))
$ws.Range("O44").Value = 32988.6
$ws.Range("N45").Copy($ws.Range("O45"))
$ws.Range("O45").Value = 30080.1
$ws.Range("N46").Copy($ws.Range("O46"))
$ws.Range("O46").Value = 35496.6
$ws.Range("N47").Copy($ws.Range("O47"))
$ws.Range("O47").Value = 23005.9
$ws.Range("N48").Copy($ws.Range("O48"))
$ws.Range("O48").Value = 20251.400000000001
$ws.Range("N49").Copy($ws.Range("O49"))
$ws.Range("O49").Value = 23657
$ws.Range("N50").Copy($ws.Range("O50"))
$ws.Range("O50").Value = 39599.9
$ws.Range("N51").Copy($ws.Range("O51"))
$ws.Range("O51").Value = 34405.599999999999
$ws.Range("N52").Copy($ws.Range("O52"))
$ws.Range("O52").Value = 42453.4
$ws.Range("N53").Copy($ws.Range("O53"))
$ws.Range("O53").Value = 24206.2
$ws.Range("N54").Copy($ws.Range("O54"))
$ws.Range("O54").Value = 24561.3
$ws.Range("N55").Copy($ws.Range("O55"))
$ws.Range("O55").Value = 23074.5
$ws.Range("N56").Copy($ws.Range("O56"))
$ws.Range("O56").Value = 20180.2
$ws.Range("N57").Copy($ws.Range("O57"))
$ws.Range("O57").Value = 20337.5
$ws.Range("N58").Copy($ws.Range("O58"))
$ws.Range("O58").Value = 19575.5
$ws.Range("N59").Copy($ws.Range("O59"))
$ws.Range("O59").Value = 20104.099999999999
$ws.Range("N60").Copy($ws.Range("O60"))
$ws.Range("O60").Value = 20456.099999999999
$ws.Range("N61").Copy($ws.Range("O61"))
$ws.Range("O61").Value = 19611.099999999999
$ws.Range("N62").Copy($ws.Range("O62"))
$ws.Range("O62").Value = 17344.5
$ws.Range("N63").Copy($ws.Range("O63"))
$ws.Range("O63").Value = 25032.3
$ws.Range("N64").Copy($ws.Range("O64"))
$ws.Range("O64").Value = 7721.8

$ws.Range("P3").Select()
